# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The error table rolls forward by one quarter: the oldest observation
# (previously in row 11) drops out, every remaining row shifts down by
# one, and a brand-new observation is inserted at the top (row 2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 2 (the newly observed quarter)
$newRow = @{ B = -0.02314597604078636; C = 0.3579920056255013; D = 0.1782699060034266; E = 0.4222202103209018; F = 0.4363822494547141; G = 15 }

# Capture current (pre-edit) values for data rows 2-11, columns B-G.
# Note: use .Value2 (not .Value) for reliable read-back of numeric values.
$cols = @("B","C","D","E","F","G")
$oldValues = @{}
for ($r = 2; $r -le 11; $r++) {
    $oldValues[$r] = @{}
    foreach ($c in $cols) {
        $oldValues[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Shift rows 2-10 down into rows 3-11 (row 11's old data is discarded)
for ($r = 10; $r -ge 2; $r--) {
    foreach ($c in $cols) {
        $ws.Range("$c$($r+1)").Value2 = $oldValues[$r][$c]
    }
}

# Write the brand-new observation into row 2
foreach ($c in $cols) {
    $ws.Range("$c" + "2").Value2 = $newRow[$c]
}
